# Applies the "Añadidas las subrutinas para calcular el nband y de pin_jointed" edit.
#
# Summary of content changes (from the OOXML diff):
#   Sheet "Elementos" (A=EA, B=Conectiv_i, C=Conectiv_j):
#     A9, A10, A11, A12, A13, A17, A18, A19, A20 : 1 -> 2
#   Sheet "Nodos" (A=X, B=Y, C=Z, D=FX, E=FY, F=Rest_x, G=Rest_y, H=Rest_z):
#     B4, B6, B8, B10, B12 : 3 -> 4
#
# Plus the view-state deltas captured in the diff: the active/selected sheet
# moves from "Nodos" to "Elementos", with "Nodos" left with a B4 selection
# and "Elementos" ending with an A21 selection.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Data edits
# ---------------------------------------------------------------------------

$wsElementos = $wb.Worksheets.Item("Elementos")
foreach ($r in 9, 10, 11, 12, 13, 17, 18, 19, 20) {
    $wsElementos.Cells.Item($r, 1).Value = 2
}

$wsNodos = $wb.Worksheets.Item("Nodos")
foreach ($r in 4, 6, 8, 10, 12) {
    $wsNodos.Cells.Item($r, 2).Value = 4
}

# ---------------------------------------------------------------------------
# 2) View-state edits
# ---------------------------------------------------------------------------

# "Nodos" keeps its own last selection (B4) even after we move away from it.
$wsNodos.Activate()
$wsNodos.Range("B4").Select()

# "Elementos" becomes the active/selected sheet, ending on A21.
$wsElementos.Activate()
$wsElementos.Range("A21").Select()
